$d = $word.ActiveDocument

# The data-quality table's 4th data row ("Access" / "" / "Drop column") is
# immediately followed by the "Wheelchair" row. Insert a new row between
# them for the Postcode/Extra characters/Strip-to-5-chars entry.
$t = $d.Tables.Item(1)

$wheelchairRow = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    if ($t.Cell($i, 1).Range.Text.TrimEnd([char]13, [char]7) -eq "Wheelchair") {
        $wheelchairRow = $t.Rows.Item($i)
        break
    }
}

$newRow = $t.Rows.Add($wheelchairRow)
$newRow.Cells.Item(1).Range.Text = "Postcode"
$newRow.Cells.Item(2).Range.Text = "Extra characters"
$newRow.Cells.Item(3).Range.Text = "Strip to just first 5 chars"
